$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 581, shifting rows 581:608 down to 582:609
$ws.Rows.Item(581).Insert()

# Populate the newly inserted row 581 with the new record
$ws.Cells.Item(581, 1).Value = 9
$ws.Cells.Item(581, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(581, 3).Value = "Metropolitana"
$ws.Cells.Item(581, 4).Value = 45147
$ws.Cells.Item(581, 5).Value = 13
$ws.Cells.Item(581, 6).Value = 100112039
$ws.Cells.Item(581, 7).Value = "Ciboulette"
$ws.Cells.Item(581, 8).Value = "Sin especificar"
$ws.Cells.Item(581, 9).Value = "Primera"
$ws.Cells.Item(581, 10).Value = 340
$ws.Cells.Item(581, 11).Value = 1200
$ws.Cells.Item(581, 12).Value = 1500
$ws.Cells.Item(581, 13).Value = 1350
$ws.Cells.Item(581, 14).Value = "`$/docena de atados"
$ws.Cells.Item(581, 15).Value = "Región Metropolitana"
$ws.Cells.Item(581, 16).Value = 450
$ws.Cells.Item(581, 17).Value = 3
$ws.Cells.Item(581, 18).Value = "Hortaliza"
